# Simulated Wild Card round and logged it
$wb = $excel.ActiveWorkbook

# OFF sheet (sheet1): row 2 - update Short Att, Short Comp, Deep Att, Deep Comp
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 501
$wsOff.Range("C2").Value = 354
$wsOff.Range("D2").Value = 130
$wsOff.Range("E2").Value = 62

# DEF sheet (sheet2): row 2 - update Short Att, Short Comp, Deep Att, Deep Comp, Short Int, Deep Int
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 600
$wsDef.Range("C2").Value = 419
$wsDef.Range("D2").Value = 129
$wsDef.Range("E2").Value = 58
$wsDef.Range("F2").Value = 13
$wsDef.Range("G2").Value = 9
